$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add new row data below the existing table (row 15, column B)
$ws.Range("B15").Value = "Cambios agregados"

# Apply the same formatting used by B10 ("17 años") onto the new cell
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Match the row height used by B10's row
$ws.Rows.Item(15).RowHeight = $ws.Rows.Item(10).RowHeight

# Restore the cell value (PasteSpecial formats only, shouldn't touch it, but make sure)
$ws.Range("B15").Value = "Cambios agregados"

# Expand the Excel Table (ListObject) to include the new row
$table = $ws.ListObjects.Item("Tabla1")
$table.Resize($ws.Range("A6:B15"))

# Scroll the view and select the new cell, matching the saved workbook view
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("B15").Select()
